$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New rows 17 and 18: "Improve VGA control code" work on 31.3.2020 ---
# Seed formatting + base content for the two new rows by duplicating row 16
# (same column styles: s=2 date, s=5 time, s=3 duration formula, s=6 task/notes).
$ws.Range("A16:F16").Copy($ws.Range("A17"))
$ws.Range("A16:F16").Copy($ws.Range("A18"))

# Row 17: 31.3.2020, 11:00-13:00, VGA Control / Improve RTL
$ws.Range("A17").Value = "31.3.2020"
$ws.Range("B17").Value = 0.45833333333333331
$ws.Range("C17").Value = 0.54166666666666663
$ws.Range("E17").Value = "VGA Control"
$ws.Range("F17").Value = "Improve RTL"

# Row 18: 31.3.2020, 14:00-14:45, VGA Control / Improve TB
$ws.Range("A18").Value = "31.3.2020"
$ws.Range("B18").Value = 0.58333333333333337
$ws.Range("C18").Value = 0.61458333333333337
$ws.Range("E18").Value = "VGA Control"
$ws.Range("F18").Value = "Improve TB"

# Duration column: extend the C-B "time spent" formula down through the new rows
$ws.Range("D17:D18").Formula = "=C17-B17"

# --- Row 38: stray formatted (time-format) cell left below the table ---
$ws.Range("B7").Copy($ws.Range("D38"))
$ws.Range("D38").ClearContents()

# --- Selection / view bookkeeping to match where the author ended up ---
$ws.Range("D38").Select()
$excel.ActiveWindow.ScrollRow = 2
$excel.ActiveWindow.ScrollColumn = 1
